$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.3025718484990989
$ws.Range("C2").Value = 0.04309411454018175
$ws.Range("D2").Value = 0.1772371813364657
$ws.Range("E2").Value = 0.1615492402676253
$ws.Range("F2").Value = 1.558199454538709
$ws.Range("I2").Value = 0.8049193884203945
$ws.Range("J2").Value = 0.1895558923542069
$ws.Range("K2").Value = 0.3220167570321166
$ws.Range("M2").Value = 0.2139809486866682
$ws.Range("O2").Value = 3.843357082501726

$ws.Range("B3").Value = 0.2701544156444413
$ws.Range("C3").Value = 0.03790341381315443
$ws.Range("D3").Value = 0.1736252473569238
$ws.Range("E3").Value = 0.1603666926146907
$ws.Range("F3").Value = 1.564409928655657
$ws.Range("I3").Value = 0.8123006490207594
$ws.Range("J3").Value = 0.1895573463846745
$ws.Range("K3").Value = 0.2864403919674601
$ws.Range("M3").Value = 0.2030191540003656
$ws.Range("O3").Value = 3.868572058196321

$ws.Range("B4").Value = 0.2502584296813666
$ws.Range("C4").Value = 0.03471053536303259
$ws.Range("D4").Value = 0.1714742075439233
$ws.Range("E4").Value = 0.159710537953675
$ws.Range("F4").Value = 1.568984973308076
$ws.Range("I4").Value = 0.8172196170061099
$ws.Range("J4").Value = 0.1896458857224843
$ws.Range("K4").Value = 0.2645940953208878
$ws.Range("M4").Value = 0.1963630088338277
$ws.Range("O4").Value = 3.885972563945444

$ws.Range("B5").Value = 0.2421532912149189
$ws.Range("C5").Value = 0.03340802738127024
$ws.Range("D5").Value = 0.1706145094831868
$ws.Range("E5").Value = 0.15946077961285
$ws.Range("F5").Value = 1.571041033467566
$ws.Range("I5").Value = 0.8193214236659756
$ws.Range("J5").Value = 0.1897040295657249
$ws.Range("K5").Value = 0.2556915205097141
$ws.Range("M5").Value = 0.1936694803239938
$ws.Range("O5").Value = 3.893545944850473

$ws.Range("B6").Value = 0.2408076135645558
$ws.Range("C6").Value = 0.03319166527978723
$ws.Range("D6").Value = 0.1704727788734886
$ws.Range("E6").Value = 0.1594203735475297
$ws.Range("F6").Value = 1.571394022390791
$ws.Range("I6").Value = 0.8196763033330079
$ws.Range("J6").Value = 0.189715017582138
$ws.Range("K6").Value = 0.2542132687093925
$ws.Range("M6").Value = 0.1932233692416254
$ws.Range("O6").Value = 3.894832646777573

$ws.Range("B7").Value = 0.2501491095899837
$ws.Range("C7").Value = 0.03469297480204148
$ws.Range("D7").Value = 0.1714625449290992
$ws.Range("E7").Value = 0.1597070981829205
$ws.Range("F7").Value = 1.569011925722023
$ws.Range("I7").Value = 0.8172475687777983
$ws.Range("J7").Value = 0.1896465805048848
$ws.Range("K7").Value = 0.2644740314967606
$ws.Range("M7").Value = 0.1963266061750062
$ws.Range("O7").Value = 3.88607274727454

$ws.Range("B8").Value = 0.2913928612525751
$ws.Range("C8").Value = 0.04130559595199657
$ws.Range("D8").Value = 0.1759779981410503
$ws.Range("E8").Value = 0.1611270088996726
$ws.Range("F8").Value = 1.560182795438941
$ws.Range("I8").Value = 0.807384165058032
$ws.Range("J8").Value = 0.1895382165538564
$ws.Range("K8").Value = 0.3097508076795066
$ws.Range("M8").Value = 0.2101859793158098
$ws.Range("O8").Value = 3.851653118463162

$ws.Range("B9").Value = 0.3723202437822124
$ws.Range("C9").Value = 0.05422488829779581
$ws.Range("D9").Value = 0.1853584434145574
$ws.Range("E9").Value = 0.16446469926478
$ws.Range("F9").Value = 1.548908488752915
$ws.Range("I9").Value = 0.791110994519844
$ws.Range("J9").Value = 0.1900203076769671
$ws.Range("K9").Value = 0.3985011947384294
$ws.Range("M9").Value = 0.2379481681597184
$ws.Range("O9").Value = 3.799373690857834

$ws.Range("B10").Value = 0.4317881612395809
$ws.Range("C10").Value = 0.06368527213962238
$ws.Range("D10").Value = 0.192566744385033
$ws.Range("E10").Value = 0.1672524102339707
$ws.Range("F10").Value = 1.544302173874485
$ws.Range("I10").Value = 0.7810253816209709
$ws.Range("J10").Value = 0.1907971141277329
$ws.Range("K10").Value = 0.4636639504626316
$ws.Range("M10").Value = 0.2586943497033118
$ws.Range("O10").Value = 3.770237344603004

$ws.Range("B11").Value = 0.4588402955390336
$ws.Range("C11").Value = 0.06798182917664519
$ws.Range("D11").Value = 0.1959139221464739
$ws.Range("E11").Value = 0.168593118261537
$ws.Range("F11").Value = 1.543004106692621
$ws.Range("I11").Value = 0.7768431767387192
$ws.Range("J11").Value = 0.1912421149287837
$ws.Range("K11").Value = 0.4932953036928893
$ws.Range("M11").Value = 0.2682068844891035
$ws.Range("O11").Value = 3.758995532142421

$ws.Range("B12").Value = 0.4690837789742943
$ws.Range("C12").Value = 0.0696077601963907
$ws.Range("D12").Value = 0.1971911171033298
$ws.Range("E12").Value = 0.1691112057528201
$ws.Range("F12").Value = 1.542627120182715
$ws.Range("I12").Value = 0.7753178332257171
$ws.Range("J12").Value = 0.1914237823454243
$ws.Range("K12").Value = 0.5045138142900782
$ws.Range("M12").Value = 0.2718196598879885
$ws.Range("O12").Value = 3.755027872163481

$ws.Range("B13").Value = 0.4668776958376384
$ws.Range("C13").Value = 0.06925763612919411
$ws.Range("D13").Value = 0.1969156208699161
$ws.Range("E13").Value = 0.1689991649281062
$ws.Range("F13").Value = 1.542703217291745
$ws.Range("I13").Value = 0.775643747138087
$ws.Range("J13").Value = 0.1913840722627924
$ws.Range("K13").Value = 0.5020978162128245
$ws.Range("M13").Value = 0.2710411164637279
$ws.Range("O13").Value = 3.75586950920777

$ws.Range("B14").Value = 0.459683047740981
$ws.Range("C14").Value = 0.06811561759695905
$ws.Range("D14").Value = 0.1960188040823709
$ws.Range("E14").Value = 0.1686355336444549
$ws.Range("F14").Value = 1.542970796324482
$ws.Range("I14").Value = 0.776716515971426
$ws.Range("J14").Value = 0.1912567972837209
$ws.Range("K14").Value = 0.4942183055426028
$ws.Range("M14").Value = 0.2685038986967996
$ws.Range("O14").Value = 3.758663311063344

$ws.Range("B15").Value = 0.4552760269953069
$ws.Range("C15").Value = 0.06741595522838395
$ws.Range("D15").Value = 0.1954707370651363
$ws.Range("E15").Value = 0.1684141508148009
$ws.Range("F15").Value = 1.543149612842484
$ws.Range("I15").Value = 0.7773812193544316
$ws.Range("J15").Value = 0.1911805502742467
$ws.Range("K15").Value = 0.4893915687348738
$ws.Range("M15").Value = 0.2669511521220187
$ws.Range("O15").Value = 3.760412279374691

$ws.Range("B16").Value = 0.4300201907129235
$ws.Range("C16").Value = 0.06340433458467487
$ws.Range("D16").Value = 0.1923493600175732
$ws.Range("E16").Value = 0.167166248214059
$ws.Range("F16").Value = 1.54440304241021
$ws.Range("I16").Value = 0.7813068653346846
$ws.Range("J16").Value = 0.1907698742955048
$ws.Range("K16").Value = 0.4617271902075402
$ws.Range("M16").Value = 0.2580741752195479
$ws.Range("O16").Value = 3.771012517713302

$ws.Range("B17").Value = 0.4145261268048444
$ws.Range("C17").Value = 0.06094148581738068
$ws.Range("D17").Value = 0.1904518640100292
$ws.Range("E17").Value = 0.1664192544226317
$ws.Range("F17").Value = 1.545376137046418
$ws.Range("I17").Value = 0.7838190651000652
$ws.Range("J17").Value = 0.1905413912375948
$ws.Range("K17").Value = 0.4447526253518959
$ws.Range("M17").Value = 0.2526475032521915
$ws.Range("O17").Value = 3.778030846356529

$ws.Range("B18").Value = 0.4056143728963377
$ws.Range("C18").Value = 0.05952426400597233
$ws.Range("D18").Value = 0.1893668889529323
$ws.Range("E18").Value = 0.1659964360963606
$ws.Range("F18").Value = 1.546010893294408
$ws.Range("I18").Value = 0.7853022100650016
$ws.Range("J18").Value = 0.1904185994958283
$ws.Range("K18").Value = 0.4349882559562275
$ws.Range("M18").Value = 0.2495332956209992
$ws.Range("O18").Value = 3.782257026545665

$ws.Range("B19").Value = 0.4025970268520496
$ws.Range("C19").Value = 0.05904430653923498
$ws.Range("D19").Value = 0.1890006397549655
$ws.Range("E19").Value = 0.165854451952395
$ws.Range("F19").Value = 1.546238705085315
$ws.Range("I19").Value = 0.7858109369011608
$ws.Range("J19").Value = 0.1903785064171757
$ws.Range("K19").Value = 0.4316820459834503
$ws.Range("M19").Value = 0.2484800988654854
$ws.Range("O19").Value = 3.783720472713838

$ws.Range("B20").Value = 0.4161754981256252
$ws.Range("C20").Value = 0.06120372872889845
$ws.Range("D20").Value = 0.190653192546236
$ws.Range("E20").Value = 0.166498066343685
$ws.Range("F20").Value = 1.54526478194667
$ws.Range("I20").Value = 0.7835476840941773
$ws.Range("J20").Value = 0.1905648210746662
$ws.Range("K20").Value = 0.4465597095291685
$ws.Range("M20").Value = 0.2532244514132813
$ws.Range("O20").Value = 3.777264128829756

$ws.Range("B21").Value = 0.4617963102508327
$ws.Range("C21").Value = 0.06845108606265171
$ws.Range("D21").Value = 0.1962819586817801
$ws.Range("E21").Value = 0.1687420593255311
$ws.Range("F21").Value = 1.542889093437822
$ws.Range("I21").Value = 0.7763998335899878
$ws.Range("J21").Value = 0.1912938241653919
$ws.Range("K21").Value = 0.4965327742660577
$ws.Range("M21").Value = 0.2692488553708969
$ws.Range("O21").Value = 3.75783485045767

$ws.Range("B22").Value = 0.4916085532684633
$ws.Range("C22").Value = 0.07318131348817758
$ws.Range("D22").Value = 0.2000171255104419
$ws.Range("E22").Value = 0.1702691814416433
$ws.Range("F22").Value = 1.542004161047615
$ws.Range("I22").Value = 0.7720685126266389
$ws.Range("J22").Value = 0.191846939167327
$ws.Range("K22").Value = 0.5291797355980918
$ws.Range("M22").Value = 0.2797833345281262
$ws.Range("O22").Value = 3.746823356915513

$ws.Range("B23").Value = 0.4756977101975508
$ws.Range("C23").Value = 0.07065730607051535
$ws.Range("D23").Value = 0.1980184649351315
$ws.Range("E23").Value = 0.1694486025635626
$ws.Range("F23").Value = 1.542415402224506
$ws.Range("I23").Value = 0.7743490859011075
$ws.Range("J23").Value = 0.191544722881325
$ws.Range("K23").Value = 0.5117568498988305
$ws.Range("M23").Value = 0.2741553135865331
$ws.Range("O23").Value = 3.752546074829638

$ws.Range("B24").Value = 0.4154298298023491
$ws.Range("C24").Value = 0.06108517274492442
$ws.Range("D24").Value = 0.1905621534663311
$ws.Range("E24").Value = 0.1664624147934042
$ws.Range("F24").Value = 1.545314890990184
$ws.Range("I24").Value = 0.7836702544386469
$ws.Range("J24").Value = 0.1905542017612092
$ws.Range("K24").Value = 0.4457427437553463
$ws.Range("M24").Value = 0.2529635955101099
$ws.Range("O24").Value = 3.777610166088152

$ws.Range("B25").Value = 0.3504240201937137
$ws.Range("C25").Value = 0.05073524187081091
$ws.Range("D25").Value = 0.1827649196524277
$ws.Range("E25").Value = 0.1635027187863116
$ws.Range("F25").Value = 1.551312420968422
$ws.Range("I25").Value = 0.7951848321158437
$ws.Range("J25").Value = 0.1898156137720903
$ws.Range("K25").Value = 0.374497867609108
$ws.Range("M25").Value = 0.2303759349418186
$ws.Range("O25").Value = 3.811887905216537
